# Generate Report for Handback
#
# For each localized-language sheet (zh-cn, de-de):
#   - "Status" column gets set to "Handed back: in sync with en-US"
#   - "Latest Target File" column gets a hyperlink to the same source doc
#     as the "Source File Name" column (the handed-back target mirrors the source)
#   - "Latest Handback File" column is filled in with the same xlf file name
#     that was produced at handoff time ("Latest Handoff File" column)
#   - "Latest Handback DateTime" column is stamped with the handback time
# The Overview sheet's per-language status columns (zh-cn / de-de) get the
# same updated status text, and a couple of columns are widened so the
# longer status / new hyperlink text aren't clipped.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# Timestamp used for this handback run per language sheet (zh-cn finished
# a few seconds before de-de).
$handbackTimes = @{
    "zh-cn" = "2016-09-01 23:07:15"
    "de-de" = "2016-09-01 23:07:22"
}

# ---- Overview sheet: refresh the per-language status cells ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---- Per-language sheets ----
$langSheets = "zh-cn", "de-de"

foreach ($langName in $langSheets) {
    $ws = $wb.Worksheets.Item($langName)
    $handbackTime = $handbackTimes[$langName]

    # Widen the Status / Latest Target File / Latest Handback File columns
    # (C, I, J) so the new, longer content is fully visible.
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664

    # Collect the existing "Source File Name" hyperlinks (column A) before
    # mutating the sheet, keyed by row number.
    $sourceLinks = @{}
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Column -eq 1) {
            $sourceLinks[$hl.Range.Row] = $hl
        }
    }

    $rows = 2, 3
    foreach ($r in $rows) {
        # Status
        $ws.Cells.Item($r, 3).Value = $statusText

        # Latest Handback File (J) mirrors the Latest Handoff File (G)
        $handoffFile = $ws.Cells.Item($r, 7).Value2
        $ws.Cells.Item($r, 10).Value = $handoffFile

        # Latest Handback DateTime (K)
        $ws.Cells.Item($r, 11).Value = $handbackTime

        # Latest Target File (I) - hyperlink to the same doc as column A
        $srcLink = $sourceLinks[$r]
        $displayName = $srcLink.TextToDisplay
        $targetUrl = $srcLink.Address
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 9), $targetUrl, "", "", $displayName)
    }
}
